## Commit: delete the "iPhone11Pro" device column from the DeviceList sheet.
## The DeviceList sheet lists several Apple devices across columns B:F
## (iPhone11, iPhone11Pro, iPhone13ProMax, iPhone7, iPhoneSE). Column C
## (iPhone11Pro / iOS 13.0.0 / UDID S2325486GUID) is removed entirely,
## shifting columns D:F left into C:E.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DeviceList")

# Remove the whole column C ("APPLE_iPhone11Pro_iOS_13.0.0_3f7e0") - this
# shifts D:F left to C:E, updates the used-range dimension, and drops the
# now-unused shared strings ("APPLE_iPhone11Pro_iOS_13.0.0_3f7e0", "13.0.0").
[void]$ws.Range("C1").EntireColumn.Delete()

# The conditional formatting that highlighted "Working"/"Not Working" text
# in row 2 was scoped to B2:F2; re-scope it to the new B2:E2 extent.
$fcs = $ws.Range("B2:E2").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    [void]$fcs.Item($i).ModifyAppliesToRange($ws.Range("B2:E2"))
}

# Leave the cursor where the author ended up after the edit.
[void]$ws.Range("D16").Select()
